# "fix cross sheet ranges"
#
# The second worksheet (originally named "Лист2", holding some stray
# COLUMN()-array-formula scratch data in F18:N28) is repurposed: it is
# renamed to "Sheet2" and now holds a plain copy of Sheet1's "Number Sold"
# column (D2:D17) in A2:A17 - so formulas that used to reach across sheets
# with positional ranges can instead use this same-sheet column.
#
# Sheet1 gains two new aggregate array formulas (F18, A20) and its
# selection/active-tab state is handed off to the new Sheet2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rebuild the second sheet -------------------------------------------
# Drop all the old COLUMN()-array-formula scratch content.
$ws2.Cells.Clear()

# Rename "Лист2" -> "Sheet2" (keeps the same r:id / physical sheet part).
$ws2.Name = "Sheet2"

# Copy Sheet1's D2:D17 ("Number Sold") values + cell formatting into the
# new Sheet2!A2:A17.
$ws1.Range("D2:D17").Copy()
$ws2.Range("A2:A17").PasteSpecial(-4104)   # xlPasteValues
$ws1.Range("D2:D17").Copy()
$ws2.Range("A2:A17").PasteSpecial(-4122)   # xlPasteFormats

for ($r = 2; $r -le 17; $r++) {
    $ws2.Rows.Item($r).RowHeight = 17.25
}

# --- New aggregate formulas on Sheet1 -----------------------------------
$ws1.Range("F18").FormulaArray = "=COUNTIFS(C2,""Coupe"",D2,8)"
$ws1.Range("A20").FormulaArray = "=SUM(A2:A17)"

# --- Selection / active tab handoff from Sheet1 to Sheet2 --------------
$null = $ws1.Range("D2:D17").Select()

$null = $ws2.Activate()
$null = $ws2.Range("A2:A17").Select()
